$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 291, shifting rows 291..308 down to 292..309
$ws.Rows.Item(291).Insert()

# Fill in the new row 291 with the new record's data.
# Columns A,B,C,E,F,G,H,I,J,K,R stay the same as the rest of the Pina/Ecuador block,
# so copy them from the row that is now 292 (the old row 291).
$ws.Cells.Item(291, 1).Value = $ws.Cells.Item(292, 1).Value2   # A Mercado ID
$ws.Cells.Item(291, 2).Value = $ws.Cells.Item(292, 2).Value2   # B Mercado
$ws.Cells.Item(291, 3).Value = $ws.Cells.Item(292, 3).Value2   # C Region
$ws.Cells.Item(291, 4).Value = 44826                           # D Fecha
$ws.Cells.Item(291, 5).Value = $ws.Cells.Item(292, 5).Value2   # E Codreg
$ws.Cells.Item(291, 6).Value = $ws.Cells.Item(292, 6).Value2   # F Tipo
$ws.Cells.Item(291, 7).Value = $ws.Cells.Item(292, 7).Value2   # G Producto ID
$ws.Cells.Item(291, 8).Value = $ws.Cells.Item(292, 8).Value2   # H Producto
$ws.Cells.Item(291, 9).Value = $ws.Cells.Item(292, 9).Value2   # I Categoria ID
$ws.Cells.Item(291, 10).Value = $ws.Cells.Item(292, 10).Value2 # J Categoria
$ws.Cells.Item(291, 11).Value = $ws.Cells.Item(292, 11).Value2 # K Variedad
$ws.Cells.Item(291, 12).Value = "Primera"                      # L Calidad
$ws.Cells.Item(291, 13).Value = 110                            # M Volumen
$ws.Cells.Item(291, 14).Value = 23000                          # N Precio minimo
$ws.Cells.Item(291, 15).Value = 23500                          # O Precio maximo
$ws.Cells.Item(291, 16).Value = 23250                          # P Precio promedio ponderado
$ws.Cells.Item(291, 17).Value = "$/caja 12 unidades"           # Q Unidad de comercializacion
$ws.Cells.Item(291, 18).Value = $ws.Cells.Item(292, 18).Value2 # R Origen
$ws.Cells.Item(291, 19).Value = 1938                           # S Precio $/Kg
$ws.Cells.Item(291, 20).Value = 12                             # T Kg / unidad

$ws.Cells.Item(291, 4).NumberFormat = $ws.Cells.Item(292, 4).NumberFormat
